# Updates the "Estado de Cuenta" detail table (rows 16-37) on Hoja1.
# The previous account-statement periods are replaced: for each worker the
# same set of periods is kept but listed newest-first (descending) instead
# of oldest-first (ascending), matching the refreshed source database.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$data = @(
    @(16, "CC", "1143390056", "JAIME LUIS ORTEGA GARCIA", "2106", 24578, 877803),
    @(17, "CC", "1143390056", "JAIME LUIS ORTEGA GARCIA", "2105", 35112, 877803),
    @(18, "CC", "1143390056", "JAIME LUIS ORTEGA GARCIA", "2104", 35112, 877803),
    @(19, "CC", "1143390056", "JAIME LUIS ORTEGA GARCIA", "2103", 35112, 877803),
    @(20, "CC", "1143390056", "JAIME LUIS ORTEGA GARCIA", "2102", 35112, 877803),
    @(21, "CC", "1143390056", "JAIME LUIS ORTEGA GARCIA", "2101", 35112, 877803),
    @(22, "CC", "1143390056", "JAIME LUIS ORTEGA GARCIA", "2012", 35112, 877803),
    @(23, "CC", "1143390056", "JAIME LUIS ORTEGA GARCIA", "2011", 35112, 877803),
    @(24, "CC", "1143390056", "JAIME LUIS ORTEGA GARCIA", "2010", 35112, 877803),
    @(25, "CC", "1143390056", "JAIME LUIS ORTEGA GARCIA", "2009", 35112, 877803),
    @(26, "CC", "1143390056", "JAIME LUIS ORTEGA GARCIA", "2008", 35112, 877803),
    @(27, "CC", "1143390056", "JAIME LUIS ORTEGA GARCIA", "2007", 35112, 877803),
    @(28, "CC", "1143390056", "JAIME LUIS ORTEGA GARCIA", "2006", 35112, 877803),
    @(29, "CC", "1143390056", "JAIME LUIS ORTEGA GARCIA", "2005", 35112, 877803),
    @(30, "CC", "1143390056", "JAIME LUIS ORTEGA GARCIA", "2003", 35112, 877803),
    @(31, "CC", "1047408531", "CARLOS CABALLERO MONTES", "2106", 25200, 900000),
    @(32, "CC", "1047408531", "CARLOS CABALLERO MONTES", "2105", 36000, 900000),
    @(33, "CC", "1047408531", "CARLOS CABALLERO MONTES", "2104", 36000, 900000),
    @(34, "CC", "1047408531", "CARLOS CABALLERO MONTES", "2103", 36000, 900000),
    @(35, "CC", "1047408531", "CARLOS CABALLERO MONTES", "2102", 36000, 900000),
    @(36, "CC", "1047408531", "CARLOS CABALLERO MONTES", "2101", 36000, 900000),
    @(37, "CC", "1047408531", "CARLOS CABALLERO MONTES", "2012", 36000, 900000)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]   # B: Tipo Doc Trabajador
    $ws.Cells.Item($r, 3).Value = $row[2]   # C: N Doc Trabajador
    $ws.Cells.Item($r, 4).Value = $row[3]   # D: Nombre Trabajador
    $ws.Cells.Item($r, 5).Value = $row[4]   # E: Periodo Mora
    $ws.Cells.Item($r, 6).Value = $row[5]   # F: Valor Mora
    $ws.Cells.Item($r, 7).Value = $row[6]   # G: Salario Basico
}
